# Insert a new data row at row 846 (2026/02/24, 火, 10:00, rank 38),
# pushing the existing rows 846..887 down to 847..888.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(846).Insert()

# Column A holds a text date like "2026/02/24" (not a real date value) in
# every other row of this sheet, so force the cell to text before writing
# it — otherwise Excel's COM layer auto-converts the "yyyy/mm/dd"-looking
# string into a date serial number. Clear the formatting back off
# afterwards so the cell matches its siblings (no explicit style/index).
$ws.Range("A846").NumberFormat = "@"
$ws.Range("A846").Value = "2026/02/24"
$ws.Range("A846").Style = "Normal"

$ws.Range("B846").Value = "火"
$ws.Range("C846").Value = 10
$ws.Range("D846").Value = 38
